$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 335, pushing existing rows 335:424 down to 336:425
$ws.Rows(335).Insert()

# Populate the newly inserted row 335 with the new data record
$ws.Range("A335").Value = 8
$ws.Range("B335").Value = "Terminal La Palmera de La Serena"
$ws.Range("C335").Value = "Coquimbo"
$ws.Range("D335").Value = 44722
$ws.Range("E335").Value = 4
$ws.Range("F335").Value = 100114001
$ws.Range("G335").Value = "Papa"
$ws.Range("H335").Value = "Asterix"
$ws.Range("I335").Value = "1a (guarda)"
$ws.Range("J335").Value = 2540
$ws.Range("K335").Value = 8500
$ws.Range("L335").Value = 9000
$ws.Range("M335").Value = 8750
$ws.Range("N335").Value = "`$/saco 25 kilos"
$ws.Range("O335").Value = "Región de La Araucanía"
$ws.Range("P335").Value = 350
$ws.Range("Q335").Value = 25
$ws.Range("R335").Value = "Hortaliza"
